$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 1.27141243209172
$ws.Range("D2").Value = 0.2168599573729835

$ws.Range("C3").Value = -0.3323401897318211
$ws.Range("D3").Value = 0.7427781762418328

$ws.Range("C4").Value = -0.07202454091174415
$ws.Range("D4").Value = 0.9432331616470808

$ws.Range("C5").Value = 0.2074131642902062
$ws.Range("D5").Value = 0.8375951108322517

$ws.Range("C6").Value = -1.992463347437849
$ws.Range("D6").Value = 0.05887219492345674

$ws.Range("C7").Value = -1.248632426695003
$ws.Range("D7").Value = 0.2249239914046868

$ws.Range("C8").Value = -0.5902159746409876
$ws.Range("D8").Value = 0.5610600912928025

$ws.Range("C9").Value = 0.2803589043019701
$ws.Range("D9").Value = 0.781821593292739

$ws.Range("C10").Value = 0.5033970510100696
$ws.Range("D10").Value = 0.6196883415105927

$ws.Range("C11").Value = 0.3508812520361659
$ws.Range("D11").Value = 0.7290151754155669
